# "csv may not have any spaces" - replace spaces in header/label text with
# underscores, drop the helper formulas in column B (replace with their
# computed literal values), and update the selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Degree:BComHons_Information_Systems_Management"

# B1 used to be "=SUM(B3:B6)/COUNTA(B3:B6)"; replace with its literal value
# and strip the bold style it had (it becomes a plain/default cell).
$ws.Range("B1").Style = "Normal"
$ws.Range("B1").Value = 53

# --- Column headings (row 2) ----------------------------------------------
$ws.Range("A2").Value = "Module_Name"
$ws.Range("B2").Value = "FinalModuleMark"

$ws.Range("C2").Value = "Mark1(%)"
$ws.Range("D2").Value = "Weight1(%)"
$ws.Range("E2").Value = "Mark2(%)"
$ws.Range("F2").Value = "Weight2(%)"
$ws.Range("G2").Value = "Mark3(%)"
$ws.Range("H2").Value = "Weight3(%)"

# Make the whole header row bold with the default (non-percent) number
# format, matching columns A2/B2 which were already bold.
$ws.Range("C2:H2").Style = "Normal"
$ws.Range("C2:H2").Font.Bold = $true

# --- Data rows (3-6): replace formulas in column B with their literal
#     values, and underscore the module names in column A -----------------
$ws.Range("A3").Value = "Advanced_Information_Systems_Theory_and_Practice"
$ws.Range("B3").Value = 75

$ws.Range("A4").Value = "Computing_in_Information_Systems"
$ws.Range("B4").Value = 50

$ws.Range("A5").Value = "Information_and_Knowledge_in_Organisations"
$ws.Range("B5").Value = 86

$ws.Range("A6").Value = "Research_Assignment:_Information_Systems_Management"
$ws.Range("B6").Value = 0

# --- A few cells in otherwise-empty columns get touched (but remain
#     valueless) in the source edit; materialize them as empty cells too.
$ws.Range("D1").Style = "Normal"
$ws.Range("F1").Style = "Normal"
$ws.Range("H1").Style = "Normal"
$ws.Range("H4").Style = "Normal"
$ws.Range("H5").Style = "Normal"
$ws.Range("F6").Style = "Normal"
$ws.Range("H6").Style = "Normal"

# --- Selection --------------------------------------------------------
$ws.Range("B4").Select() | Out-Null
